# Adding to E2E tests - Manager user stories + acceptance criteria scenarios.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Move the "Reimbursement System" block down (old rows 19-21 -> new 21,23,25),
#    working bottom-up so we don't clobber data we still need to read.
# ---------------------------------------------------------------------------

# old row21 -> new row25
$ws.Range("A21:C21").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)
$ws.Range("A25:C25").PasteSpecial(-4163)

# old row20 -> new row23
$ws.Range("A20:C20").Copy()
$ws.Range("A23:C23").PasteSpecial(-4122)
$ws.Range("A23:C23").PasteSpecial(-4163)

# old row19 -> new row21
$ws.Range("A19:C19").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A21:C21").PasteSpecial(-4163)

# Clear the now-stale cells left behind at rows 19 & 20 (content only, for now)
$ws.Range("A19:C20").ClearContents()

# New blank separator rows 19, 22, 24 (style like row10 - A=B=C same plain style)
$ws.Range("A10:C10").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Range("A24:C24").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Manager section (rows 12-18)
# ---------------------------------------------------------------------------

# Row 18 ("to log out") takes the style that row 17 used to have (bottom of box)
$ws.Range("A17:C17").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Range("B18").Value = "to log out"
$ws.Range("C18").Value = "my information will not remain on the computer."

# Row 17 becomes a blank "mid" row (style like row 13/14/15/16)
$ws.Range("A13:C13").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A17:C17").ClearContents()

# Row 14 gets the new "approve or deny" story, row 13 & 15 are cleared (blank, same style)
$ws.Range("B14").Value = "to approve or deny reimbursement requests"
$ws.Range("C14").Value = "if they are legitimate or not."
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()

# Row 16 keeps its existing story text (unchanged)
# (B16 / C16 already contain "to view reimbursement statistics" / "I can keep track of employee activities.")

Write-Host "checkpoint2"
